# Update "想去人数" (F column) figures across sheets, matching the
# regenerated data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (sheet1) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 219
$ws1.Range("F8").Value  = 59
$ws1.Range("F9").Value  = 256
$ws1.Range("F13").Value = 96
$ws1.Range("F14").Value = 2105
$ws1.Range("F15").Value = 52
$ws1.Range("F16").Value = 26
$ws1.Range("F17").Value = 501
$ws1.Range("F18").Value = 481
$ws1.Range("F23").Value = 1584
$ws1.Range("F24").Value = 3786
$ws1.Range("F25").Value = 27
$ws1.Range("F26").Value = 59
$ws1.Range("F28").Value = 1131
$ws1.Range("F29").Value = 114
$ws1.Range("F30").Value = 1912
$ws1.Range("F32").Value = 461
$ws1.Range("F33").Value = 70
$ws1.Range("F34").Value = 279
$ws1.Range("F35").Value = 409
$ws1.Range("F37").Value = 656
$ws1.Range("F38").Value = 434
$ws1.Range("F39").Value = 385

# ---- Sheet "演出" (sheet2) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 19

# ---- Sheet "全部类型" (sheet4) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 219
$ws4.Range("F8").Value  = 59
$ws4.Range("F9").Value  = 256
$ws4.Range("F13").Value = 96
$ws4.Range("F14").Value = 2105
$ws4.Range("F15").Value = 52
$ws4.Range("F16").Value = 19
$ws4.Range("F17").Value = 26
$ws4.Range("F18").Value = 501
$ws4.Range("F19").Value = 481
$ws4.Range("F24").Value = 1584
$ws4.Range("F25").Value = 3786
$ws4.Range("F26").Value = 27
$ws4.Range("F27").Value = 59
$ws4.Range("F29").Value = 1131
$ws4.Range("F30").Value = 114
$ws4.Range("F31").Value = 1912
$ws4.Range("F33").Value = 461
$ws4.Range("F34").Value = 70
$ws4.Range("F35").Value = 279
$ws4.Range("F36").Value = 409
$ws4.Range("F38").Value = 656
$ws4.Range("F39").Value = 434
$ws4.Range("F40").Value = 385
